$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.291738
$ws.Range("H2").Value = 18.875214
$ws.Range("I2").Value = 0.5742845621220376
$ws.Range("J2").Value = 0.5742845621220376
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 132.641248155486
$ws.Range("R2").Value = 1193.771233399374
$ws.Range("S2").Value = 0.03282645098000485
$ws.Range("T2").Value = 0.03282645098000485
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.291738
$ws.Range("H3").Value = 18.875214
$ws.Range("I3").Value = 0.5742845621220376
$ws.Range("J3").Value = 0.5742845621220376
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 1897.589116171716
$ws.Range("R3").Value = 17078.30204554544
$ws.Range("S3").Value = 0.469621003785957
$ws.Range("T3").Value = 0.469621003785957
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.291738
$ws.Range("H4").Value = 18.875214
$ws.Range("I4").Value = 0.5742845621220376
$ws.Range("J4").Value = 0.5742845621220376
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 290.270903467254
$ws.Range("R4").Value = 2612.438131205286
$ws.Range("S4").Value = 0.07183710735607576
$ws.Range("T4").Value = 0.07183710735607576
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.730603666666667
$ws.Range("H5").Value = 11.191811
$ws.Range("I5").Value = 0.3405145117553424
$ws.Range("J5").Value = 0.3405145117553424
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 78.64789136485011
$ws.Range("R5").Value = 707.831022283651
$ws.Range("S5").Value = 0.01946401429774407
$ws.Range("T5").Value = 0.01946401429774407
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.730603666666667
$ws.Range("H6").Value = 11.191811
$ws.Range("I6").Value = 0.3405145117553424
$ws.Range("J6").Value = 0.3405145117553424
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 1125.150620483078
$ws.Range("R6").Value = 10126.3555843477
$ws.Range("S6").Value = 0.2784556252449755
$ws.Range("T6").Value = 0.2784556252449755
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.730603666666667
$ws.Range("H7").Value = 11.191811
$ws.Range("I7").Value = 0.3405145117553424
$ws.Range("J7").Value = 0.3405145117553424
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 172.1123315690488
$ws.Range("R7").Value = 1549.010984121439
$ws.Range("S7").Value = 0.04259487221262284
$ws.Range("T7").Value = 0.04259487221262284
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.933443
$ws.Range("H8").Value = 2.800329
$ws.Range("I8").Value = 0.08520092612262004
$ws.Range("J8").Value = 0.08520092612262004
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 19.67867139445433
$ws.Range("R8").Value = 177.108042550089
$ws.Range("S8").Value = 0.00487013618210559
$ws.Range("T8").Value = 0.00487013618210559
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.933443
$ws.Range("H9").Value = 2.800329
$ws.Range("I9").Value = 0.08520092612262004
$ws.Range("J9").Value = 0.08520092612262004
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 281.5265475718593
$ws.Range("R9").Value = 2533.738928146734
$ws.Range("S9").Value = 0.06967302812624669
$ws.Range("T9").Value = 0.06967302812624669
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.933443
$ws.Range("H10").Value = 2.800329
$ws.Range("I10").Value = 0.08520092612262004
$ws.Range("J10").Value = 0.08520092612262004
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 43.06462585460234
$ws.Range("R10").Value = 387.5816326914211
$ws.Range("S10").Value = 0.01065776181426776
$ws.Range("T10").Value = 0.01065776181426776
